$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: read a cell's raw value (Value2) without forcing Excel's
# automatic type inference to corrupt it on the way back out.
function Get-CellValue($ws, $addr) {
    return $ws.Range($addr).Value2
}

# Helper: write a value back to a cell, forcing TEXT storage (so that
# numeric-looking strings like "15" do not silently become real numbers),
# while clearing the "quote prefix" cell style that Excel applies when a
# leading apostrophe is used, and properly blanking out $null/"" values.
function Set-CellText($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    if ($null -eq $val -or $val -eq "") {
        $cell.ClearContents()
    } else {
        $cell.Formula = "'" + $val
        $cell.Style = "Normal"
    }
}

# Helper: write a value back to a cell as a genuine NUMBER.
function Set-CellNumber($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value2 = $val
    }
}

# Columns that hold genuine numbers vs. columns that hold text
# (including numeric-looking text such as "15" in column I).
$numericCols = @("A", "B", "E", "Q", "R")
$textCols    = @("D", "F", "G", "H", "I", "J", "K", "P", "AC")

function Swap-Row($ws, $r1, $r2) {
    foreach ($col in $numericCols) {
        $a1 = "$col$r1"
        $a2 = "$col$r2"
        $v1 = Get-CellValue $ws $a1
        $v2 = Get-CellValue $ws $a2
        Set-CellNumber $ws $a1 $v2
        Set-CellNumber $ws $a2 $v1
    }
    foreach ($col in $textCols) {
        $a1 = "$col$r1"
        $a2 = "$col$r2"
        $v1 = Get-CellValue $ws $a1
        $v2 = Get-CellValue $ws $a2
        Set-CellText $ws $a1 $v2
        Set-CellText $ws $a2 $v1
    }
}

# The edit swaps the full content of row 2 <-> row 7, row 3 <-> row 6,
# and row 4 <-> row 8 (rows 1, 5, 9 are untouched).
Swap-Row $ws 2 7
Swap-Row $ws 3 6
Swap-Row $ws 4 8
